# Planting randomization update: clear stale CMS_cold (J) dates for the
# first batch, backfill the second batch's CMS_heat/chlorophyll (I/K)
# dates, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing date-formatted cell as the style donor so new cells pick
# up the same cellXf (s="1", numFmtId 14) instead of Excel minting a new
# number-format entry.
$dateDonor = $ws.Range("I2")

# --- Rows 2-6: the "CMS_cold" (J) date was recorded in error; clear it. ---
$ws.Range("J2").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("J6").Value = ""

# Row 5's "chlorophyll" (K) reading actually happened on 3/5/2021.
$ws.Range("K5").Value = "3/5/2021"

# --- Rows 7-11: backfill CMS_heat (I) / CMS_cold (J) / chlorophyll (K). ---
$dateDonor.Copy($ws.Range("I7"))
$ws.Range("I7").Value = "3/5/2021"
$dateDonor.Copy($ws.Range("J7"))
$ws.Range("J7").Value = ""
$dateDonor.Copy($ws.Range("K7"))
$ws.Range("K7").Value = "3/5/2021"

$dateDonor.Copy($ws.Range("I8"))
$ws.Range("I8").Value = "3/5/2021"

$dateDonor.Copy($ws.Range("I9"))
$ws.Range("I9").Value = "3/5/2021"
$dateDonor.Copy($ws.Range("K9"))
$ws.Range("K9").Value = "3/5/2021"

$dateDonor.Copy($ws.Range("I10"))
$ws.Range("I10").Value = "3/5/2021"
$dateDonor.Copy($ws.Range("K10"))
$ws.Range("K10").Value = "3/5/2021"

$dateDonor.Copy($ws.Range("I11"))
$ws.Range("I11").Value = "3/5/2021"
$dateDonor.Copy($ws.Range("K11"))
$ws.Range("K11").Value = "3/5/2021"

# Reflect where the author's cursor ended up after the edits.
$ws.Range("M10").Select() | Out-Null
